# Add a new row of dynamic test data to the "ContactUs" sheet and make it
# the active sheet/selection, mirroring what happens when a user appends a
# row in Excel and saves.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactUs")

$ws.Range("A7").Value = "John@123#"
$ws.Range("B7").Value = "validemail@example.com"
$ws.Range("C7").Value = "Invalid Name TC"
$ws.Range("D7").Value = "Testing invalid characters."
$ws.Range("E7").Value = "INVALID_NAME"

$ws.Activate()
$ws.Range("A7:E7").Select()
